# Append rows 206-217 (feature index 204-215) to the "월_전체승객" sheet,
# mirroring the existing A/B column pattern (A = index, B = value),
# including the same direct cell formatting used by the existing A column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing formatting of the last populated "A" cell (A205) onto
# the new A206:A217 range so the new index cells keep the same bold/border/
# centered style used throughout column A, without introducing new style
# definitions.
$ws.Range("A205").Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @(204, 0.6086956521739131),
    @(205, 0.4347826086956522),
    @(206, 0.4641304347826087),
    @(207, 0.6376811594202898),
    @(208, 0.4021739130434783),
    @(209, 0.732919254658385),
    @(210, 0.6413043478260869),
    @(211, 0.4217391304347826),
    @(212, 0.4184782608695652),
    @(213, 0.6521739130434783),
    @(214, 0.4347826086956522),
    @(215, 0.4347826086956522)
)

$startRow = 206
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
}
